$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 233, shifting the existing rows 233:245 down to 234:246.
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new weekly price record.
$ws.Range("A233").Value = 11
$ws.Range("B233").Value = "Vega Monumental Concepción"
$ws.Range("C233").Value = "Bíobío"
$ws.Range("D233").Value = 45267
$ws.Range("E233").Value = 8
$ws.Range("F233").Value = 100112021
$ws.Range("G233").Value = "Ají"
$ws.Range("H233").Value = "Inferno"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 50
$ws.Range("K233").Value = 34000
$ws.Range("L233").Value = 35000
$ws.Range("M233").Value = 34600
$ws.Range("N233").Value = "$/caja 10 kilos"
$ws.Range("O233").Value = "Región de Arica y Parinacota"
$ws.Range("P233").Value = 3460
$ws.Range("Q233").Value = 10
$ws.Range("R233").Value = "Hortaliza"
